# Update "想去人数" (want-to-go count, column F) figures that changed
# between crawler runs, as published to gh-pages (output generated at 456a3b4).
#
# Sheet "展览" (Exhibitions)
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 130
$ws1.Range("F8").Value = 928
$ws1.Range("F13").Value = 901
$ws1.Range("F15").Value = 3747
$ws1.Range("F16").Value = 1122
$ws1.Range("F18").Value = 2533
$ws1.Range("F20").Value = 1060
$ws1.Range("F21").Value = 3492
$ws1.Range("F22").Value = 728
$ws1.Range("F24").Value = 34
$ws1.Range("F25").Value = 2085
$ws1.Range("F27").Value = 812
$ws1.Range("F30").Value = 175
$ws1.Range("F34").Value = 473

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 298

# Sheet "全部类型" (All types) - combined listing, mirrors the same events
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 130
$ws4.Range("F6").Value = 928
$ws4.Range("F14").Value = 901
$ws4.Range("F16").Value = 3747
$ws4.Range("F17").Value = 1122
$ws4.Range("F20").Value = 2533
$ws4.Range("F22").Value = 1060
$ws4.Range("F23").Value = 3492
$ws4.Range("F24").Value = 728
$ws4.Range("F27").Value = 34
$ws4.Range("F28").Value = 2085
$ws4.Range("F34").Value = 812
$ws4.Range("F37").Value = 175
$ws4.Range("F44").Value = 473

$wb.Save()
